$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values to update: column D (Price) values are forced to text so that
# numeric-looking strings (e.g. "6.93") are not auto-converted to numbers by Excel,
# matching the original inline-string storage. NumberFormat is reset to "Normal"
# style afterwards so the cell keeps its original default styling.
$values = @{
    'D2' = '63.537.28'
    'E2' = '  -0.12%  '
    'D3' = '2.648.28'
    'E3' = '  -0.33%  '
    'E4' = '  +0.03%  '
    'D5' = '602.08'
    'E5' = '  +1.79%  '
    'D6' = '146.90'
    'E6' = '  +1.50%  '
    'E7' = '  +0.00%  '
    'E8' = '  +0.18%  '
    'E9' = '  +1.20%  '
    'E10' = '  -0.55%  '
    'E11' = '  +4.41%  '
    'D13' = '27.49'
    'E13' = '  +0.04%  '
    'D14' = '3.126.30'
    'E14' = '  -0.25%  '
    'D15' = '63.421.52'
    'E15' = '  -0.19%  '
    'E16' = '  +0.17%  '
    'D17' = '2.663.13'
    'E17' = '  +0.84%  '
    'E19' = '  +4.29%  '
    'E20' = '  +0.14%  '
    'D21' = '6.93'
    'E21' = '  +2.72%  '
    'E22' = '  -0.10%  '
    'E23' = '  -3.36%  '
    'E24' = '  -1.28%  '
    'E25' = '  +1.50%  '
    'D26' = '9.12'
    'E26' = '  +7.54%  '
    'D27' = '1.56'
    'E27' = '  +0.73%  '
    'D28' = '561.36'
    'E28' = '  +2.45%  '
    'E29' = '  -1.58%  '
    'E30' = '  +0.01%  '
    'D31' = '7.94'
    'E31' = '  +1.37%  '
    'D32' = '2.03'
    'E32' = '  +2.86%  '
    'E33' = '  -3.87%  '
    'D34' = '0.0₃0817'
    'E34' = '  +0.74%  '
    'D35' = '5.16'
    'E35' = '  +4.84%  '
    'D36' = '167.56'
    'E36' = '  -3.93%  '
    'E37' = '  +0.83%  '
    'E38' = '  -0.03%  '
    'D39' = '1.92'
    'E39' = '  +5.24%  '
    'E40' = '  -0.04%  '
    'E41' = '  +0.02%  '
    'D42' = '168.52'
    'E42' = '  -1.64%  '
    'D43' = '3.77'
    'E43' = '  +0.72%  '
    'E44' = '  -0.93%  '
    'E45' = '  +2.66%  '
    'E46' = '  -0.06%  '
    'E47' = '  +3.15%  '
    'E48' = '  -0.29%  '
    'D49' = '18.79'
    'E49' = '  +0.12%  '
    'E50' = '  +9.15%  '
    'E51' = '  -0.73%  '
}

foreach ($addr in $values.Keys) {
    $col = $addr.Substring(0,1)
    $range = $ws.Range($addr)
    if ($col -eq "D") {
        $range.NumberFormat = "@"
        $range.Value = $values[$addr]
        $range.Style = "Normal"
    } else {
        $range.Value = $values[$addr]
    }
}
